$d = $word.ActiveDocument

# 1. Total de citas programadas: 7 -> 6
$d.Content.Find.Execute("7", $false, $true, $false, $false, $false, $true, 1, $false, "6", 2) | Out-Null

# 2. Swap the buyers in the first two rows (08:30-08:45 and 08:45-09:00):
#    REGIONAL S.A.S <-> COLFRESH COFFEE  (use a temp placeholder to do a clean swap)
$d.Content.Find.Execute("REGIONAL S.A.S", $false, $true, $false, $false, $false, $true, 1, $false, "__TEMP_SWAP__", 2) | Out-Null
$d.Content.Find.Execute("COLFRESH COFFEE", $false, $true, $false, $false, $false, $true, 1, $false, "REGIONAL S.A.S", 2) | Out-Null
$d.Content.Find.Execute("__TEMP_SWAP__", $false, $true, $false, $false, $false, $true, 1, $false, "COLFRESH COFFEE", 2) | Out-Null

# 3. Shift the BOX BRAND and ARMANDO VELÁSQUEZ time slots forward an hour
$d.Content.Find.Execute("09:00 - 09:15", $false, $true, $false, $false, $false, $true, 1, $false, "10:00 - 10:15", 2) | Out-Null
$d.Content.Find.Execute("09:15 - 09:30", $false, $true, $false, $false, $false, $true, 1, $false, "10:15 - 10:30", 2) | Out-Null

# 4. Remove the INMERSSO BOUTIQUE row (09:30 - 09:45) entirely - it is row 6 of the table
#    (row 1 = header; rows 2-5 = 08:30, 08:45, 09:00(now 10:00), 09:15(now 10:15); row 6 = 09:30)
$t = $d.Tables(1)
$t.Rows(6).Delete()

# 5. Update the (previously) last row FIRST, before touching the row that will collide with it:
#    10:30 - 10:45 -> 11:30 - 11:45, and buyer ENCADENAMIENTOS... -> INMERSSO BOUTIQUE
$d.Content.Find.Execute("10:30 - 10:45", $false, $true, $false, $false, $false, $true, 1, $false, "11:30 - 11:45", 2) | Out-Null
$d.Content.Find.Execute("ENCADENAMIENTOS PRODUCTIVOS -  CAFE AROMAS DEL EJE / CAFÉ GRANEAO.", $false, $true, $false, $false, $false, $true, 1, $false, "INMERSSO BOUTIQUE", 2) | Out-Null

# 6. Now shift the INTERLINK2AMERICAS row's time slot: 09:45 - 10:00 -> 10:30 - 10:45
$d.Content.Find.Execute("09:45 - 10:00", $false, $true, $false, $false, $false, $true, 1, $false, "10:30 - 10:45", 2) | Out-Null

Write-Output "edit complete"
